$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(105, 8).Value = 47900
$ws.Cells.Item(105, 10).Value = 47900
$ws.Cells.Item(105, 12).Value = 47900
$ws.Cells.Item(105, 14).Value = -54888
$ws.Cells.Item(112, 8).Value = 1526.3
$ws.Cells.Item(112, 9).Value = 1065.5555
$ws.Cells.Item(112, 10).Value = 1903.2727
$ws.Cells.Item(112, 11).Value = 3196.6665
$ws.Cells.Item(112, 12).Value = 5709.8181
$ws.Cells.Item(112, 13).Value = -2088.6665
$ws.Cells.Item(112, 14).Value = -7925.8181
$ws.Cells.Item(115, 8).Value = 579
$ws.Cells.Item(115, 9).Value = 587.7778
$ws.Cells.Item(115, 10).Value = 500
$ws.Cells.Item(115, 11).Value = 1763.3334
$ws.Cells.Item(115, 12).Value = 1500
$ws.Cells.Item(115, 13).Value = -196.3334
$ws.Cells.Item(115, 14).Value = -4634
$ws.Cells.Item(118, 8).Value = 1781.1765
$ws.Cells.Item(118, 9).Value = 410
$ws.Cells.Item(118, 11).Value = 1230
$ws.Cells.Item(118, 13).Value = 427
$ws.Cells.Item(132, 8).Value = 2860.743
$ws.Cells.Item(132, 9).Value = 2433.4614
$ws.Cells.Item(132, 10).Value = 4095.111
$ws.Cells.Item(132, 11).Value = 7300.3842
$ws.Cells.Item(132, 12).Value = 12285.333
$ws.Cells.Item(132, 13).Value = -4770.3842
$ws.Cells.Item(132, 14).Value = -17345.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(29, 8).Value = 3316.5557
$ws.Cells.Item(29, 9).Value = 1524.75
$ws.Cells.Item(29, 10).Value = 4750
$ws.Cells.Item(29, 11).Value = 1524.75
$ws.Cells.Item(29, 12).Value = 4750
$ws.Cells.Item(29, 13).Value = -1216.75
$ws.Cells.Item(29, 14).Value = -5366
$ws.Cells.Item(74, 8).Value = 3683.7317
$ws.Cells.Item(74, 9).Value = 841.93335
$ws.Cells.Item(74, 10).Value = 11434.091
$ws.Cells.Item(74, 11).Value = 841.93335
$ws.Cells.Item(74, 12).Value = 11434.091
$ws.Cells.Item(74, 13).Value = 32.06664999999998
$ws.Cells.Item(74, 14).Value = -13182.091
$ws.Cells.Item(77, 8).Value = 3683.7317
$ws.Cells.Item(77, 9).Value = 841.93335
$ws.Cells.Item(77, 10).Value = 11434.091
$ws.Cells.Item(77, 11).Value = 4209.66675
$ws.Cells.Item(77, 12).Value = 57170.455
$ws.Cells.Item(77, 13).Value = 158.3332499999997
$ws.Cells.Item(77, 14).Value = -65906.455
$ws.Cells.Item(82, 8).Value = 35750
$ws.Cells.Item(82, 10).Value = 35750
$ws.Cells.Item(82, 12).Value = 35750
$ws.Cells.Item(82, 14).Value = -36472
$ws.Cells.Item(85, 8).Value = 35750
$ws.Cells.Item(85, 10).Value = 35750
$ws.Cells.Item(85, 12).Value = 35750
$ws.Cells.Item(85, 14).Value = -38246
$ws.Cells.Item(97, 8).Value = 39489
$ws.Cells.Item(97, 9).Value = 56284.055
$ws.Cells.Item(97, 10).Value = 1700.125
$ws.Cells.Item(97, 11).Value = 56284.055
$ws.Cells.Item(97, 12).Value = 1700.125
$ws.Cells.Item(97, 13).Value = -55788.055
$ws.Cells.Item(97, 14).Value = -2692.125
$ws.Cells.Item(110, 8).Value = 2836.1143
$ws.Cells.Item(110, 9).Value = 2904.5217
$ws.Cells.Item(110, 10).Value = 2705
$ws.Cells.Item(110, 11).Value = 2904.5217
$ws.Cells.Item(110, 12).Value = 2705
$ws.Cells.Item(110, 13).Value = -859.5216999999998
$ws.Cells.Item(110, 14).Value = -6795

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 79295.16
$ws.Cells.Item(20, 9).Value = 796.6667
$ws.Cells.Item(20, 10).Value = 102844.7
$ws.Cells.Item(20, 11).Value = 796.6667
$ws.Cells.Item(20, 12).Value = 102844.7
$ws.Cells.Item(20, 13).Value = -549.6667
$ws.Cells.Item(20, 14).Value = -103338.7
$ws.Cells.Item(107, 8).Value = 3283.8
$ws.Cells.Item(107, 9).Value = 2833
$ws.Cells.Item(107, 10).Value = 3960
$ws.Cells.Item(107, 11).Value = 2833
$ws.Cells.Item(107, 12).Value = 3960
$ws.Cells.Item(107, 13).Value = -913
$ws.Cells.Item(107, 14).Value = -7800

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 48300
$ws.Cells.Item(51, 9).Value = 30000
$ws.Cells.Item(51, 10).Value = 57450
$ws.Cells.Item(51, 11).Value = 30000
$ws.Cells.Item(51, 12).Value = 57450
$ws.Cells.Item(51, 13).Value = -29264
$ws.Cells.Item(51, 14).Value = -58922
$ws.Cells.Item(61, 8).Value = 48300
$ws.Cells.Item(61, 9).Value = 30000
$ws.Cells.Item(61, 10).Value = 57450
$ws.Cells.Item(61, 11).Value = 30000
$ws.Cells.Item(61, 12).Value = 57450
$ws.Cells.Item(61, 13).Value = -29652
$ws.Cells.Item(61, 14).Value = -58146

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 148.60606
$ws.Cells.Item(33, 9).Value = 34.434784
$ws.Cells.Item(33, 10).Value = 411.2
$ws.Cells.Item(33, 11).Value = 206.608704
$ws.Cells.Item(33, 12).Value = 2467.2
$ws.Cells.Item(33, 13).Value = 76.39129600000001
$ws.Cells.Item(33, 14).Value = -3033.2
$ws.Cells.Item(55, 8).Value = 3355.7144
$ws.Cells.Item(55, 9).Value = 490
$ws.Cells.Item(55, 10).Value = 3833.3333
$ws.Cells.Item(55, 11).Value = 1470
$ws.Cells.Item(55, 12).Value = 11499.9999
$ws.Cells.Item(55, 13).Value = -1293
$ws.Cells.Item(55, 14).Value = -11853.9999
$ws.Cells.Item(80, 8).Value = 128312.625
$ws.Cells.Item(80, 10).Value = 4000
$ws.Cells.Item(80, 12).Value = 12000
$ws.Cells.Item(80, 14).Value = -13872
$ws.Cells.Item(83, 8).Value = 128312.625
$ws.Cells.Item(83, 10).Value = 4000
$ws.Cells.Item(83, 12).Value = 36000
$ws.Cells.Item(83, 14).Value = -45360
$ws.Cells.Item(87, 8).Value = 2015.5
$ws.Cells.Item(87, 9).Value = 1555.1538
$ws.Cells.Item(87, 10).Value = 8000
$ws.Cells.Item(87, 11).Value = 4665.4614
$ws.Cells.Item(87, 12).Value = 24000
$ws.Cells.Item(87, 13).Value = -3417.4614
$ws.Cells.Item(87, 14).Value = -26496
$ws.Cells.Item(90, 8).Value = 2015.5
$ws.Cells.Item(90, 9).Value = 1555.1538
$ws.Cells.Item(90, 10).Value = 8000
$ws.Cells.Item(90, 11).Value = 13996.3842
$ws.Cells.Item(90, 12).Value = 72000
$ws.Cells.Item(90, 13).Value = -7756.3842
$ws.Cells.Item(90, 14).Value = -84480
$ws.Cells.Item(117, 8).Value = 666.6667
$ws.Cells.Item(117, 9).Value = 300
$ws.Cells.Item(117, 10).Value = 850
$ws.Cells.Item(117, 11).Value = 900
$ws.Cells.Item(117, 12).Value = 2550
$ws.Cells.Item(117, 13).Value = 2542
$ws.Cells.Item(117, 14).Value = -9434
$ws.Cells.Item(121, 8).Value = 9288.025
$ws.Cells.Item(121, 10).Value = 10011.194
$ws.Cells.Item(121, 12).Value = 30033.582
$ws.Cells.Item(121, 14).Value = -32653.582
$ws.Cells.Item(139, 8).Value = 459095.6
$ws.Cells.Item(139, 9).Value = 524409.25
$ws.Cells.Item(139, 10).Value = 1900
$ws.Cells.Item(139, 11).Value = 1573227.75
$ws.Cells.Item(139, 12).Value = 5700
$ws.Cells.Item(139, 13).Value = -1568087.75
$ws.Cells.Item(139, 14).Value = -15980

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 837.6
$ws.Cells.Item(31, 9).Value = 837.6
$ws.Cells.Item(31, 11).Value = 837.6
$ws.Cells.Item(31, 13).Value = -545.6
$ws.Cells.Item(37, 8).Value = 837.6
$ws.Cells.Item(37, 9).Value = 837.6
$ws.Cells.Item(37, 11).Value = 837.6
$ws.Cells.Item(37, 13).Value = -560.6
$ws.Cells.Item(80, 8).Value = 2186.2307
$ws.Cells.Item(80, 9).Value = 1987.8572
$ws.Cells.Item(80, 10).Value = 2417.6667
$ws.Cells.Item(80, 11).Value = 1987.8572
$ws.Cells.Item(80, 12).Value = 2417.6667
$ws.Cells.Item(80, 13).Value = -989.8571999999999
$ws.Cells.Item(80, 14).Value = -4413.6667
$ws.Cells.Item(83, 8).Value = 2186.2307
$ws.Cells.Item(83, 9).Value = 1987.8572
$ws.Cells.Item(83, 10).Value = 2417.6667
$ws.Cells.Item(83, 11).Value = 9939.286
$ws.Cells.Item(83, 12).Value = 12088.3335
$ws.Cells.Item(83, 13).Value = -4947.286
$ws.Cells.Item(83, 14).Value = -22072.3335
$ws.Cells.Item(132, 8).Value = 19429.08
$ws.Cells.Item(132, 9).Value = 30114.436
$ws.Cells.Item(132, 10).Value = 2065.375
$ws.Cells.Item(132, 11).Value = 90343.308
$ws.Cells.Item(132, 12).Value = 6196.125
$ws.Cells.Item(132, 13).Value = -87813.308
$ws.Cells.Item(132, 14).Value = -11256.125
